$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 51
# Leading apostrophe forces the date-looking string to stay literal text
# instead of being auto-parsed into a date serial number by Excel.
$ws.Cells.Item($row, 1).Value = "'08/17/2025"
# Re-normalize the style so the quote-prefix formatting Excel applied
# doesn't leave a stray style index on the cell (matches the plain,
# unstyled data rows elsewhere in the sheet).
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 106.9919999999984
$ws.Cells.Item($row, 3).Value = 0.09346493195753096
$ws.Cells.Item($row, 4).Value = 10
